# Binary Search 2: Square root of Integer - updated index
#
# Adds a new "Github Link" entry (column F, row 3) on the "Binary Search 2"
# sheet, mirroring the existing pattern used on "Binary Search 1": the cell
# shows a friendly description while the hyperlink points at the GitHub
# source file.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Binary Search 2")

$targetUrl  = "https://github.com/ankurnecessary/dsa/blob/main/2_binarySearch/1_square_root_of_integer.js"
$friendlyText = "dsa/1_square_root_of_integer.js at main " + [char]0x00B7 + " ankurnecessary/dsa " + [char]0x00B7 + " GitHub"

$cell = $ws.Range("F3")

# Create the hyperlink first (target only) so we can separately cache its
# display text without disturbing the cell's own value.
$ws.Hyperlinks.Add($cell, $targetUrl) | Out-Null

$hl = $ws.Hyperlinks.Item($ws.Hyperlinks.Count)
$hl.TextToDisplay = $targetUrl

# Now set the cell's visible text to the friendly description (matches the
# style used for every other Github Link cell in this workbook).
$cell.Value = $friendlyText

# Copy the formatting of an equivalent "Github Link" cell on sheet 1 so the
# new cell reuses the workbook's existing Hyperlink-ish style instead of
# minting a new one.
$srcWs = $wb.Worksheets.Item("Binary Search 1")
$srcWs.Range("F4").Copy() | Out-Null
$cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
